# Update TPM-derived values in the LR-pairs worksheet (Icam4-Itgb1)
# per the "update scripts wuth new tpm" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2564746666666666
$ws.Range("H2").Value = 0.7694239999999999
$ws.Range("I2").Value = 0.1818007399394835
$ws.Range("J2").Value = 0.1818007399394835
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 19.77150980040889
$ws.Range("R2").Value = 177.94358820368
$ws.Range("S2").Value = 0.04370156499624046
$ws.Range("T2").Value = 0.04370156499624047
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2564746666666666
$ws.Range("H3").Value = 0.7694239999999999
$ws.Range("I3").Value = 0.1818007399394835
$ws.Range("J3").Value = 0.1818007399394835
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 26.05270621505422
$ws.Range("R3").Value = 234.474355935488
$ws.Range("S3").Value = 0.05758508305529627
$ws.Range("T3").Value = 0.05758508305529626
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2564746666666666
$ws.Range("H4").Value = 0.7694239999999999
$ws.Range("I4").Value = 0.1818007399394835
$ws.Range("J4").Value = 0.1818007399394835
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 36.42627345200355
$ws.Range("R4").Value = 327.836461068032
$ws.Range("S4").Value = 0.08051409188794678
$ws.Range("T4").Value = 0.0805140918879468
$ws.Range("I5").Value = 0.7694380609030022
$ws.Range("J5").Value = 0.7694380609030022
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 83.67926427040555
$ws.Range("R5").Value = 753.11337843365
$ws.Range("S5").Value = 0.1849588040198672
$ws.Range("T5").Value = 0.1849588040198672
$ws.Range("I6").Value = 0.7694380609030022
$ws.Range("J6").Value = 0.7694380609030022
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2437182305075021
$ws.Range("T6").Value = 0.2437182305075021
$ws.Range("I7").Value = 0.7694380609030022
$ws.Range("J7").Value = 0.7694380609030022
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.3407610263756329
$ws.Range("T7").Value = 0.3407610263756329
$ws.Range("G8").Value = 0.06878966666666667
$ws.Range("I8").Value = 0.0487611991575143
$ws.Range("J8").Value = 0.0487611991575143
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 5.302962613592222
$ws.Range("R8").Value = 47.72666352233
$ws.Range("S8").Value = 0.01172129835657472
$ws.Range("T8").Value = 0.01172129835657472
$ws.Range("G9").Value = 0.06878966666666667
$ws.Range("I9").Value = 0.0487611991575143
$ws.Range("J9").Value = 0.0487611991575143
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("Q9").Value = 6.987656908147556
$ws.Range("R9").Value = 62.888912173328
$ws.Range("S9").Value = 0.01544502901526134
$ws.Range("T9").Value = 0.01544502901526133
$ws.Range("G10").Value = 0.06878966666666667
$ws.Range("I10").Value = 0.0487611991575143
$ws.Range("J10").Value = 0.0487611991575143
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("Q10").Value = 9.769975495976889
$ws.Range("R10").Value = 87.929779463792
$ws.Range("S10").Value = 0.02159487178567824
$ws.Range("T10").Value = 0.02159487178567824
